$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Requirement #4 currently reads "... Field/scroll/button for every hole
# color". The edit removes the word "every " (leaving "for hole") and,
# because the author's cursor ended up inside the word "hole" after the
# retype, the "_GoBack" bookmark (previously sitting at the very end of
# requirement #8's "** Optional Low **" line) is now anchored between
# "ho" and "le".
# --------------------------------------------------------------------

# Locate paragraph 4's sentence in the current document text.
$fullText = $d.Content.Text
$anchorIdx = $fullText.IndexOf("for every hole")

# 1) Drop a throw-away bookmark right at the run boundary that precedes
#    "for every " (i.e. immediately after "Field/scroll/button "). This
#    keeps that run boundary intact while we edit the neighbouring run,
#    so "Field/scroll/button " and "for " don't get coalesced into one
#    run by the edit below.
$boundaryRange = $d.Range($anchorIdx, $anchorIdx)
$d.Bookmarks.Add("TempBoundary", $boundaryRange) | Out-Null

# 2) Remove the word "every " (turns "for every hole" into "for hole").
$fullText = $d.Content.Text
$everyIdx = $fullText.IndexOf("for every ho")
$everyRange = $d.Range($everyIdx + 4, $everyIdx + 10)
$everyRange.Delete()

# 3) Drop the temporary bookmark now that the edit is done; this does not
#    re-merge the runs it protected.
$d.Bookmarks("TempBoundary").Delete()

# 4) Re-seat the "_GoBack" bookmark between "ho" and "le" of "hole".
#    Adding a bookmark with a name that already exists simply moves it,
#    which removes it from its old spot at the end of requirement #8.
$fullText = $d.Content.Text
$forHoleIdx = $fullText.IndexOf("for hole")
$holeIdx = $fullText.IndexOf("hole", $forHoleIdx)
$splitPos = $holeIdx + 2
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $splitRange) | Out-Null
